$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric -> force text to match original inline-string typing
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.73'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4684'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07379'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8721'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.39'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.390'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07081'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.521'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.88'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008722'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.74'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.332'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.84'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.178'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.331'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.15'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08960'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7690'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.167'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.512'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.085'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01963'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05294'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.967'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.281'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.337'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1654'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.50'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.32'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06303'

# Plain text / percentage cells (never auto-converted to numbers)
$ws.Range("D2").Value = '26.929.89'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.816.71'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("E8").Value = '  -1.21%  '
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '1.813.93'
$ws.Range("E12").Value = '  +2.61%  '
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '26.958.13'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("D24").Value = '2.040.76'
$ws.Range("E24").Value = '  +1.94%  '
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E31").Value = '  +0.81%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("E40").Value = '  +2.13%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("E47").Value = '  +2.18%  '
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("E51").Value = '  -0.30%  '
